# acwe2023.xlsx update: "Updated settings, all plots complete now."
#
# 1. Coureurs sheet (sheet1): the Saudi-Arabië (G) column now holds what used
#    to be the Bahrein (F) result, and the Australië (H) column is reset to 0
#    (race not yet scored / plot reset). Selection moves to G2:G21.
# 2. Teams sheet (sheet4): 19 new team rows (ids 3-21) are appended with the
#    corresponding participant names, chassis/motor picks and points.
#    Selection moves to F23 (last filled cell).

$wb = $excel.ActiveWorkbook

# --- Coureurs: copy Bahrein (F) results into Saudi-Arabië (G), reset
#     Australië (H) to 0, for every data row (2-21) ---
$ws1 = $wb.Worksheets.Item("Coureurs")

for ($r = 2; $r -le 21; $r++) {
    $fVal = $ws1.Cells.Item($r, 6).Value2
    $ws1.Cells.Item($r, 7).Value = $fVal
    $ws1.Cells.Item($r, 8).Value = 0
}

$ws1.Range("G2:G21").Select()

# --- Teams: append the 19 missing rows (id 3 .. 21) ---
$ws4 = $wb.Worksheets.Item("Teams")

$teamRows = @(
    @(3,  "Casper",    5,  13, 5, 0),
    @(4,  "Raymond",   1,  13, 3, 1),
    @(5,  "Niels",     6,  13, 2, 0),
    @(6,  "Erik",      3,  8,  6, 0),
    @(7,  "Grietje",   9,  16, 0, 2),
    @(8,  "Arjan T.",  16, 5,  4, 0),
    @(9,  "Charlotte", 1,  16, 3, 1),
    @(10, "Michiel",   1,  8,  4, 1),
    @(11, "Emily",     3,  9,  4, 2),
    @(12, "Rodi",      9,  5,  3, 2),
    @(13, "Arjan Z.",  9,  5,  4, 2),
    @(14, "Hans",      5,  7,  4, 1),
    @(15, "Kitty",     5,  6,  3, 1),
    @(16, "Stan",      9,  4,  5, 2),
    @(17, "Erwin",     2,  11, 3, 0),
    @(18, "Mees",      4,  6,  3, 2),
    @(19, "Jordi",     2,  8,  4, 1),
    @(20, "Wietse",    2,  7,  7, 3),
    @(21, "Isolde",    2,  16, 4, 1)
)

$row = 5
foreach ($team in $teamRows) {
    $ws4.Cells.Item($row, 1).Value = $team[0]
    $ws4.Cells.Item($row, 2).Value = $team[1]
    $ws4.Cells.Item($row, 3).Value = $team[2]
    $ws4.Cells.Item($row, 4).Value = $team[3]
    $ws4.Cells.Item($row, 5).Value = $team[4]
    $ws4.Cells.Item($row, 6).Value = $team[5]
    $row++
}

$ws4.Range("F23").Select()

# Coureurs stays the tab-selected sheet, so re-activate it last and restore
# its own selection.
$ws1.Activate()
$ws1.Range("G2:G21").Select()
